$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 436
$ws.Range("I4").Value = 436
$ws.Range("K4").Value = 436
$ws.Range("M4").Value = -322

$ws.Range("H6").Value = 85.57143000000001
$ws.Range("I6").Value = 83.166664
$ws.Range("K6").Value = 249.499992
$ws.Range("M6").Value = -137.499992

$ws.Range("H9").Value = 140
$ws.Range("I9").Value = 175.8
$ws.Range("J9").Value = 95.25
$ws.Range("K9").Value = 175.8
$ws.Range("L9").Value = 95.25
$ws.Range("M9").Value = -6.800000000000011
$ws.Range("N9").Value = -433.25

$ws.Range("H41").Value = 359.66666
$ws.Range("I41").Value = 229.7
$ws.Range("J41").Value = 619.6
$ws.Range("K41").Value = 229.7
$ws.Range("L41").Value = 619.6
$ws.Range("M41").Value = 210.3
$ws.Range("N41").Value = -1499.6

$ws.Range("H64").Value = 3966.9333
$ws.Range("I64").Value = 3891.3333
$ws.Range("K64").Value = 3891.3333
$ws.Range("M64").Value = -3643.3333

$ws.Range("H67").Value = 3966.9333
$ws.Range("I67").Value = 3891.3333
$ws.Range("K67").Value = 3891.3333
$ws.Range("M67").Value = -3033.3333

$ws.Range("H100").Value = 1697.0476
$ws.Range("I100").Value = 1466.8125
$ws.Range("J100").Value = 2433.8
$ws.Range("K100").Value = 1466.8125
$ws.Range("L100").Value = 2433.8
$ws.Range("M100").Value = -925.8125
$ws.Range("N100").Value = -3515.8

$ws.Range("H128").Value = 150000
$ws.Range("J128").Value = 150000
$ws.Range("L128").Value = 150000
$ws.Range("N128").Value = -159960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3157.2856
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 3157.2856
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3157.2856
$ws.Range("N5").Value = -3381.2856
$ws.Range("M5").ClearContents()

$ws.Range("H45").Value = 4043.125
$ws.Range("I45").Value = 3279
$ws.Range("K45").Value = 3279
$ws.Range("M45").Value = -2902

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3157.2856
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3157.2856
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3157.2856
$ws.Range("N4").Value = -3387.2856
$ws.Range("M4").ClearContents()

$ws.Range("H99").Value = 8820.357
$ws.Range("I99").Value = 11629.3
$ws.Range("K99").Value = 11629.3
$ws.Range("M99").Value = -10131.3

$ws.Range("H134").Value = 3411.6086
$ws.Range("I134").Value = 3411.6086
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10234.8258
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7699.825800000001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 94.583336
$ws.Range("I7").Value = 46.833332
$ws.Range("K7").Value = 46.833332
$ws.Range("M7").Value = 66.166668

$ws.Range("H13").Value = 1800
$ws.Range("I13").Value = 1800
$ws.Range("K13").Value = 1800
$ws.Range("M13").Value = -1661

$ws.Range("H16").Value = 16441.3
$ws.Range("I16").Value = 32125
$ws.Range("J16").Value = 5985.5
$ws.Range("K16").Value = 32125
$ws.Range("L16").Value = 5985.5
$ws.Range("M16").Value = -31838
$ws.Range("N16").Value = -6559.5

$ws.Range("H70").Value = 37500
$ws.Range("J70").Value = 37500
$ws.Range("L70").Value = 37500
$ws.Range("N70").Value = -38130

$ws.Range("H73").Value = 37500
$ws.Range("J73").Value = 37500
$ws.Range("L73").Value = 37500
$ws.Range("N73").Value = -39684

$ws.Range("H113").Value = 16441.3
$ws.Range("I113").Value = 32125
$ws.Range("J113").Value = 5985.5
$ws.Range("K113").Value = 32125
$ws.Range("L113").Value = 5985.5
$ws.Range("M113").Value = -29955
$ws.Range("N113").Value = -10325.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 534.7895
$ws.Range("J12").Value = 267.33334
$ws.Range("L12").Value = 802.0000200000001
$ws.Range("N12").Value = -1148.00002

$ws.Range("H19").Value = 1207.25
$ws.Range("I19").Value = 551.6667
$ws.Range("J19").Value = 1600.6
$ws.Range("K19").Value = 1655.0001
$ws.Range("L19").Value = 4801.799999999999
$ws.Range("M19").Value = -1481.0001
$ws.Range("N19").Value = -5149.799999999999

$ws.Range("H131").Value = 2280.25
$ws.Range("I131").Value = 823.875
$ws.Range("K131").Value = 2471.625
$ws.Range("M131").Value = 2568.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 75
$ws.Range("I13").Value = 29.166666
$ws.Range("J13").Value = 166.66667
$ws.Range("K13").Value = 29.166666
$ws.Range("L13").Value = 166.66667
$ws.Range("M13").Value = 109.833334
$ws.Range("N13").Value = -444.66667

$ws.Range("H75").Value = 65000
$ws.Range("J75").Value = 65000
$ws.Range("L75").Value = 65000
$ws.Range("N75").Value = -66748

$ws.Range("H78").Value = 65000
$ws.Range("J78").Value = 65000
$ws.Range("L78").Value = 195000
$ws.Range("N78").Value = -203736

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H113").Value = 1622.6
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 16355.625
$ws.Range("I132").Value = 16355.625
$ws.Range("K132").Value = 49066.875
$ws.Range("M132").Value = -46536.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3251.0952
$ws.Range("J22").Value = 3527.0557
$ws.Range("L22").Value = 3527.0557
$ws.Range("N22").Value = -4117.0557

$ws.Range("H27").Value = 3251.0952
$ws.Range("J27").Value = 3527.0557
$ws.Range("L27").Value = 3527.0557
$ws.Range("N27").Value = -3741.0557

$ws.Range("H60").Value = 29999.5
$ws.Range("I60").Value = 29999.5
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 29999.5
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -29490.5
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H96").Value = 1680.1428
$ws.Range("J96").Value = 885.5
$ws.Range("L96").Value = 885.5
$ws.Range("N96").Value = -3631.5
